# Edit script applying the "Ya funciona el algoritmo pero falta revisar" commit.
# Changes:
#  1. class_teachers!B26: reassign teacher (Xim H. -> Yolanda M.) and drop the
#     bottom border so it matches the style used a few rows below (B29).
#  2. teacher_hours (rows 2-32): update the weekly-availability grid (B:F "X"
#     markers) and the computed hour totals in column G - this is the output
#     of re-running the genetic timetabling algorithm.
#  3. Selections / active sheet are updated on all three affected sheets to
#     match where the author was last working (teacher_hours ends up active).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("course_hours")
$ws2 = $wb.Worksheets.Item("class_teachers")
$ws3 = $wb.Worksheets.Item("teacher_hours")

# --- class_teachers: row 26 teacher reassignment --------------------------
$ws2.Range("B26").Value = "Yolanda M."
$ws2.Range("B26").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> xlLineStyleNone (matches B29's style)

# --- teacher_hours: refreshed availability grid + hour totals -------------
$ws3.Range("C2").Value = $null
$ws3.Range("D2").Value = "X"
$ws3.Range("G2").Value = 8
$ws3.Range("C3").Value = "X"
$ws3.Range("G3").Value = 2
$ws3.Range("D4").Value = $null
$ws3.Range("G4").Value = 15
$ws3.Range("C5").Value = $null
$ws3.Range("G5").Value = 2
$ws3.Range("B6").Value = $null
$ws3.Range("F6").Value = "X"
$ws3.Range("G6").Value = 10
$ws3.Range("B7").Value = "X"
$ws3.Range("C7").Value = $null
$ws3.Range("F7").Value = "X"
$ws3.Range("G7").Value = 2
$ws3.Range("B8").Value = $null
$ws3.Range("E8").Value = "X"
$ws3.Range("G8").Value = 4
$ws3.Range("B9").Value = "X"
$ws3.Range("D9").Value = "X"
$ws3.Range("G9").Value = 2
$ws3.Range("G10").Value = 2
$ws3.Range("B11").Value = $null
$ws3.Range("G11").Value = 7
$ws3.Range("E12").Value = $null
$ws3.Range("C13").Value = $null
$ws3.Range("D14").Value = "X"
$ws3.Range("G14").Value = 6
$ws3.Range("G15").Value = 19
$ws3.Range("F16").Value = $null
$ws3.Range("G16").Value = 24
$ws3.Range("F17").Value = $null
$ws3.Range("G17").Value = 24
$ws3.Range("C18").Value = $null
$ws3.Range("G18").Value = 16
$ws3.Range("G19").Value = 4
$ws3.Range("C20").Value = "X"
$ws3.Range("D20").Value = "X"
$ws3.Range("G20").Value = 4
$ws3.Range("G21").Value = 6
$ws3.Range("E22").Value = $null
$ws3.Range("F22").Value = $null
$ws3.Range("C23").Value = $null
$ws3.Range("D23").Value = $null
$ws3.Range("G23").Value = 9
$ws3.Range("C24").Value = $null
$ws3.Range("G24").Value = 30
$ws3.Range("C25").Value = $null
$ws3.Range("G25").Value = 30
$ws3.Range("F26").Value = $null
$ws3.Range("G26").Value = 30
$ws3.Range("D27").Value = $null
$ws3.Range("G27").Value = 30
$ws3.Range("G28").Value = 30
$ws3.Range("B29").Value = $null
$ws3.Range("G29").Value = 30
$ws3.Range("C30").Value = $null
$ws3.Range("E30").Value = $null
$ws3.Range("G30").Value = 30
$ws3.Range("C31").Value = "X"
$ws3.Range("D31").Value = $null
$ws3.Range("G31").Value = 8
$ws3.Range("G32").Value = 30

# --- Selections / active sheet, matching the author's final view state ----
$ws1.Range("B5").Select()
$ws2.Range("A10").Select()
$ws3.Activate()
$ws3.Range("H21").Select()
